$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml) - column F is "想去人数" (want-to-go count)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 472
$ws1.Range("F4").Value = 7919
$ws1.Range("F5").Value = 94
$ws1.Range("F10").Value = 461
$ws1.Range("F12").Value = 19
$ws1.Range("F13").Value = 448
$ws1.Range("F14").Value = 67
$ws1.Range("F17").Value = 5814
$ws1.Range("F18").Value = 174
$ws1.Range("F19").Value = 252
$ws1.Range("F20").Value = 1670
$ws1.Range("F21").Value = 233
$ws1.Range("F22").Value = 370

# Sheet "全部类型" (sheet4.xml) - same updates, rows shifted by the extra
# "演出" entry present in this combined sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 472
$ws4.Range("F4").Value = 7919
$ws4.Range("F5").Value = 94
$ws4.Range("F10").Value = 461
$ws4.Range("F12").Value = 19
$ws4.Range("F13").Value = 448
$ws4.Range("F14").Value = 67
$ws4.Range("F18").Value = 5814
$ws4.Range("F20").Value = 174
$ws4.Range("F21").Value = 252
$ws4.Range("F22").Value = 1670
$ws4.Range("F23").Value = 233
$ws4.Range("F24").Value = 370
